$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 2999.5
$ws.Cells.Item(76, 9).Value = 2999.5
$ws.Cells.Item(76, 11).Value = 2999.5
$ws.Cells.Item(76, 13).Value = -2684.5
$ws.Cells.Item(79, 8).Value = 2999.5
$ws.Cells.Item(79, 9).Value = 2999.5
$ws.Cells.Item(79, 11).Value = 2999.5
$ws.Cells.Item(79, 13).Value = -1907.5
$ws.Cells.Item(113, 8).Value = 1932.7142
$ws.Cells.Item(113, 9).Value = 1691
$ws.Cells.Item(113, 11).Value = 1691
$ws.Cells.Item(113, 13).Value = 1563
$ws.Cells.Item(116, 8).Value = 11158.474
$ws.Cells.Item(116, 9).Value = 13140.733
$ws.Cells.Item(116, 11).Value = 13140.733
$ws.Cells.Item(116, 13).Value = -9698.733
$ws.Cells.Item(127, 8).Value = 1096.25
$ws.Cells.Item(127, 9).Value = 824.2857
$ws.Cells.Item(127, 11).Value = 2472.8571
$ws.Cells.Item(127, 13).Value = 2487.1429
$ws.Cells.Item(137, 8).Value = 86098.8
$ws.Cells.Item(137, 9).Value = 138783.17
$ws.Cells.Item(137, 11).Value = 416349.51
$ws.Cells.Item(137, 13).Value = -413799.51
$ws.Cells.Item(138, 8).Value = 3252.675
$ws.Cells.Item(138, 10).Value = 3704.8276
$ws.Cells.Item(138, 12).Value = 11114.4828
$ws.Cells.Item(138, 14).Value = -21394.4828

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2817
$ws.Cells.Item(2, 9).Value = 1397.6
$ws.Cells.Item(2, 10).Value = 3999.8333
$ws.Cells.Item(2, 11).Value = 1397.6
$ws.Cells.Item(2, 12).Value = 3999.8333
$ws.Cells.Item(2, 13).Value = -1284.6
$ws.Cells.Item(2, 14).Value = -4225.8333
$ws.Cells.Item(26, 8).Value = 16779.6
$ws.Cells.Item(26, 9).Value = 9724.5
$ws.Cells.Item(26, 11).Value = 9724.5
$ws.Cells.Item(26, 13).Value = -9394.5
$ws.Cells.Item(32, 8).Value = 3351.756
$ws.Cells.Item(32, 9).Value = 2565.081
$ws.Cells.Item(32, 11).Value = 2565.081
$ws.Cells.Item(32, 13).Value = -2278.081
$ws.Cells.Item(61, 8).Value = 3893.4
$ws.Cells.Item(61, 9).Value = 3714.889
$ws.Cells.Item(61, 11).Value = 3714.889
$ws.Cells.Item(61, 13).Value = -3502.889
$ws.Cells.Item(74, 8).Value = 1425.8889
$ws.Cells.Item(74, 9).Value = 1389.04
$ws.Cells.Item(74, 11).Value = 1389.04
$ws.Cells.Item(74, 13).Value = -515.04
$ws.Cells.Item(77, 8).Value = 1425.8889
$ws.Cells.Item(77, 9).Value = 1389.04
$ws.Cells.Item(77, 11).Value = 6945.2
$ws.Cells.Item(77, 13).Value = -2577.2
$ws.Cells.Item(116, 8).Value = 2817
$ws.Cells.Item(116, 9).Value = 1397.6
$ws.Cells.Item(116, 10).Value = 3999.8333
$ws.Cells.Item(116, 11).Value = 1397.6
$ws.Cells.Item(116, 12).Value = 3999.8333
$ws.Cells.Item(116, 13).Value = 896.4000000000001
$ws.Cells.Item(116, 14).Value = -8587.8333
$ws.Cells.Item(122, 8).Value = 5271.2085
$ws.Cells.Item(122, 9).Value = 3634.1333
$ws.Cells.Item(122, 11).Value = 10902.3999
$ws.Cells.Item(122, 13).Value = -8452.3999
$ws.Cells.Item(132, 8).Value = 288610.06
$ws.Cells.Item(132, 9).Value = 325251.66
$ws.Cells.Item(132, 10).Value = 4637.75
$ws.Cells.Item(132, 11).Value = 975754.98
$ws.Cells.Item(132, 12).Value = 13913.25
$ws.Cells.Item(132, 13).Value = -973224.98
$ws.Cells.Item(132, 14).Value = -18973.25
$ws.Cells.Item(136, 8).Value = 3893.4
$ws.Cells.Item(136, 9).Value = 3714.889
$ws.Cells.Item(136, 11).Value = 11144.667
$ws.Cells.Item(136, 13).Value = -8594.667000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2817
$ws.Cells.Item(3, 9).Value = 1397.6
$ws.Cells.Item(3, 10).Value = 3999.8333
$ws.Cells.Item(3, 11).Value = 1397.6
$ws.Cells.Item(3, 12).Value = 3999.8333
$ws.Cells.Item(3, 13).Value = -1283.6
$ws.Cells.Item(3, 14).Value = -4227.8333
$ws.Cells.Item(60, 8).Value = 106662.664
$ws.Cells.Item(60, 10).Value = 106662.664
$ws.Cells.Item(60, 12).Value = 106662.664
$ws.Cells.Item(60, 14).Value = -107860.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(8, 8).Value = 10833
$ws.Cells.Item(8, 10).Value = 10833
$ws.Cells.Item(8, 12).Value = 10833
$ws.Cells.Item(8, 14).Value = -11113
$ws.Cells.Item(31, 8).Value = 3574.5356
$ws.Cells.Item(31, 9).Value = 2340.5557
$ws.Cells.Item(31, 11).Value = 2340.5557
$ws.Cells.Item(31, 13).Value = -2045.5557
$ws.Cells.Item(34, 8).Value = 3574.5356
$ws.Cells.Item(34, 9).Value = 2340.5557
$ws.Cells.Item(34, 11).Value = 2340.5557
$ws.Cells.Item(34, 13).Value = -2138.5557
$ws.Cells.Item(58, 8).Value = 3466.2632
$ws.Cells.Item(58, 9).Value = 3195.4348
$ws.Cells.Item(58, 10).Value = 3881.5334
$ws.Cells.Item(58, 11).Value = 3195.4348
$ws.Cells.Item(58, 12).Value = 3881.5334
$ws.Cells.Item(58, 13).Value = -2992.4348
$ws.Cells.Item(58, 14).Value = -4287.5334
$ws.Cells.Item(105, 8).Value = 2863.1667
$ws.Cells.Item(105, 9).Value = 2419.75
$ws.Cells.Item(105, 11).Value = 2419.75
$ws.Cells.Item(105, 13).Value = -672.75
$ws.Cells.Item(122, 8).Value = 1456.4166
$ws.Cells.Item(122, 9).Value = 1467.7
$ws.Cells.Item(122, 11).Value = 4403.1
$ws.Cells.Item(122, 13).Value = -1953.1
$ws.Cells.Item(132, 8).Value = 5245.1113
$ws.Cells.Item(132, 9).Value = 4867.8335
$ws.Cells.Item(132, 10).Value = 5999.6665
$ws.Cells.Item(132, 11).Value = 14603.5005
$ws.Cells.Item(132, 12).Value = 17998.9995
$ws.Cells.Item(132, 13).Value = -12073.5005
$ws.Cells.Item(132, 14).Value = -23058.9995
$ws.Cells.Item(134, 8).Value = 3706.55
$ws.Cells.Item(134, 9).Value = 2941.6428
$ws.Cells.Item(134, 11).Value = 8824.928400000001
$ws.Cells.Item(134, 13).Value = -6289.928400000001
$ws.Cells.Item(136, 8).Value = 3466.2632
$ws.Cells.Item(136, 9).Value = 3195.4348
$ws.Cells.Item(136, 10).Value = 3881.5334
$ws.Cells.Item(136, 11).Value = 9586.304400000001
$ws.Cells.Item(136, 12).Value = 11644.6002
$ws.Cells.Item(136, 13).Value = -7036.304400000001
$ws.Cells.Item(136, 14).Value = -16744.6002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(44, 8).Value = 2239.25
$ws.Cells.Item(44, 9).Value = 1130.5714
$ws.Cells.Item(44, 10).Value = 10000
$ws.Cells.Item(44, 11).Value = 3391.7142
$ws.Cells.Item(44, 12).Value = 30000
$ws.Cells.Item(44, 13).Value = -2993.7142
$ws.Cells.Item(44, 14).Value = -30796
$ws.Cells.Item(64, 8).Value = 1500
$ws.Cells.Item(64, 9).Value = 1500
$ws.Cells.Item(64, 11).Value = 4500
$ws.Cells.Item(64, 13).Value = -4230
$ws.Cells.Item(67, 8).Value = 1500
$ws.Cells.Item(67, 9).Value = 1500
$ws.Cells.Item(67, 11).Value = 4500
$ws.Cells.Item(67, 13).Value = -3564
$ws.Cells.Item(103, 8).Value = 326.5
$ws.Cells.Item(103, 9).Value = 341.5
$ws.Cells.Item(103, 10).Value = 281.5
$ws.Cells.Item(103, 11).Value = 1024.5
$ws.Cells.Item(103, 12).Value = 844.5
$ws.Cells.Item(103, 13).Value = -145.5
$ws.Cells.Item(103, 14).Value = -2602.5
$ws.Cells.Item(109, 8).Value = 2398
$ws.Cells.Item(109, 9).Value = 836.8
$ws.Cells.Item(109, 10).Value = 5000
$ws.Cells.Item(109, 11).Value = 2510.4
$ws.Cells.Item(109, 12).Value = 15000
$ws.Cells.Item(109, 13).Value = -1470.4
$ws.Cells.Item(109, 14).Value = -17080
$ws.Cells.Item(119, 8).Value = 3005.4
$ws.Cells.Item(119, 9).Value = 2756.75
$ws.Cells.Item(119, 10).Value = 4000
$ws.Cells.Item(119, 11).Value = 8270.25
$ws.Cells.Item(119, 12).Value = 12000
$ws.Cells.Item(119, 13).Value = -3432.25
$ws.Cells.Item(119, 14).Value = -21676

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2917.4167
$ws.Cells.Item(102, 9).Value = 2920.9
$ws.Cells.Item(102, 10).Value = 2900
$ws.Cells.Item(102, 11).Value = 2920.9
$ws.Cells.Item(102, 12).Value = 2900
$ws.Cells.Item(102, 13).Value = -1298.9
$ws.Cells.Item(102, 14).Value = -6144
$ws.Cells.Item(122, 8).Value = 2124
$ws.Cells.Item(122, 9).Value = 1872
$ws.Cells.Item(122, 11).Value = 5616
$ws.Cells.Item(122, 13).Value = -3166
$ws.Cells.Item(126, 8).Value = 3757
$ws.Cells.Item(126, 9).Value = 3086.1538
$ws.Cells.Item(126, 10).Value = 4629.1
$ws.Cells.Item(126, 11).Value = 9258.4614
$ws.Cells.Item(126, 12).Value = 13887.3
$ws.Cells.Item(126, 13).Value = -6788.4614
$ws.Cells.Item(126, 14).Value = -18827.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1836.2354
$ws.Cells.Item(16, 9).Value = 1614.2667
$ws.Cells.Item(16, 10).Value = 3501
$ws.Cells.Item(16, 11).Value = 1614.2667
$ws.Cells.Item(16, 12).Value = 3501
$ws.Cells.Item(16, 13).Value = -1444.2667
$ws.Cells.Item(16, 14).Value = -3841
$ws.Cells.Item(61, 8).Value = 5338.2
$ws.Cells.Item(61, 9).Value = 1537
$ws.Cells.Item(61, 11).Value = 1537
$ws.Cells.Item(61, 13).Value = -1335
$ws.Cells.Item(113, 8).Value = 5338.2
$ws.Cells.Item(113, 9).Value = 1537
$ws.Cells.Item(113, 11).Value = 1537
$ws.Cells.Item(113, 13).Value = 633
$ws.Cells.Item(122, 8).Value = 91666
$ws.Cells.Item(122, 9).Value = 69999.336
$ws.Cells.Item(122, 11).Value = 209998.008
$ws.Cells.Item(122, 13).Value = -207548.008
$ws.Cells.Item(132, 8).Value = 2003540.4
$ws.Cells.Item(132, 9).Value = 3336667.2
$ws.Cells.Item(132, 10).Value = 3850
$ws.Cells.Item(132, 11).Value = 10010001.6
$ws.Cells.Item(132, 12).Value = 11550
$ws.Cells.Item(132, 13).Value = -10007471.6
$ws.Cells.Item(132, 14).Value = -16610

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 10748.9
$ws.Cells.Item(41, 10).Value = 10436.25
$ws.Cells.Item(41, 12).Value = 10436.25
$ws.Cells.Item(41, 14).Value = -11216.25
$ws.Cells.Item(45, 8).Value = 31374.75
$ws.Cells.Item(45, 9).Value = 4500
$ws.Cells.Item(45, 10).Value = 40333
$ws.Cells.Item(45, 11).Value = 4500
$ws.Cells.Item(45, 12).Value = 40333
$ws.Cells.Item(45, 13).Value = -4009
$ws.Cells.Item(45, 14).Value = -41315
$ws.Cells.Item(74, 8).Value = 29624.2
$ws.Cells.Item(74, 10).Value = 31655.25
$ws.Cells.Item(74, 12).Value = 31655.25
$ws.Cells.Item(74, 14).Value = -33527.25
$ws.Cells.Item(77, 8).Value = 29624.2
$ws.Cells.Item(77, 10).Value = 31655.25
$ws.Cells.Item(77, 12).Value = 94965.75
$ws.Cells.Item(77, 14).Value = -104325.75
$ws.Cells.Item(107, 8).Value = 335.22726
$ws.Cells.Item(107, 10).Value = 411.0909
$ws.Cells.Item(107, 12).Value = 1233.2727
$ws.Cells.Item(107, 14).Value = -5073.2727
$ws.Cells.Item(122, 8).Value = 5840.846
$ws.Cells.Item(122, 9).Value = 1908.5555
$ws.Cells.Item(122, 11).Value = 5725.666499999999
$ws.Cells.Item(122, 13).Value = -3275.666499999999
$ws.Cells.Item(132, 8).Value = 21595.055
$ws.Cells.Item(132, 9).Value = 25871.592
$ws.Cells.Item(132, 11).Value = 77614.776
$ws.Cells.Item(132, 13).Value = -75084.776
